$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells for Wins / Losses / Ties, matching the style of the existing header row (row 1)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Fill in team record data (Wins=77, Losses=85, Ties=0) for every data row
for ($r = 2; $r -le 53; $r++) {
    $ws.Cells.Item($r, 30).Value = 77   # AD
    $ws.Cells.Item($r, 31).Value = 85   # AE
    $ws.Cells.Item($r, 32).Value = 0    # AF
}
